$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.95%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.70%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.250"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-7.86%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08137"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.353"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.68%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.799"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-11.69%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9554"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.81%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1856"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.42%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09375"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04621"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.98%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.462"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-28.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.12%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005935"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.15%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.380"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.36%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.518"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.04%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.38%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1367"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.34%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2712"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.16%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04189"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.19%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001257"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004315"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-8.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.46%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002995"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-19.90%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02598"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-6.83%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05469"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.07%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007834"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.39%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.85%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006600"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.30%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002131"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.03%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008692"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.69%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3441"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006989"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.84%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000756"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.85%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003494"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.06%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003550"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.49%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002118"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.85%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002017"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.85%"
